$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new range to be stored as text so numeric-looking strings
# ("0", "140.00", etc.) are not auto-converted to numbers by Excel.
$rng = $ws.Range("A8:K13")
$rng.NumberFormat = "@"

$nbsp = [char]0x00A0
$player = "Ravichandran Ashwin" + $nbsp

$data = @(
    @(" Dubai (DSC)", " October 14 2020", "Capitals won by 13 runs", "Delhi Capitals", "Rajasthan Royals", $player, "0", "0", "0", "0", "-"),
    @(" Dubai (DSC)", " October 27 2020", "Sunrisers won by 88 runs", "Delhi Capitals", "Sunrisers Hyderabad", $player, "7", "5", "1", "0", "140.00"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Delhi Capitals", "Kolkata Knight Riders", $player, "14", "13", "2", "0", "107.69"),
    @(" Dubai (DSC)", " October 31 2020", "Mumbai won by 9 wickets (with 34 balls remaining)", "Delhi Capitals", "Mumbai Indians", $player, "12", "9", "0", "1", "133.33"),
    @(" Dubai (DSC)", " September 20 2020", "Match tied (Capitals won the one-over eliminator)", "Delhi Capitals", "Kings XI Punjab", $player, "4", "6", "0", "0", "66.66"),
    @(" Sharjah", " October 09 2020", "Capitals won by 46 runs", "Delhi Capitals", "Rajasthan Royals", $player, "0", "1", "0", "0", "0.00")
)

$startRow = 8
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($c = 0; $c -lt 11; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $data[$i][$c]
    }
}
